{"js": "const replacements = [\n  [\"393\u00d73=1179\", \"297\u00d77=2079\"],\n  [\"693\u00d72=1386\", \"402\u00d78=3216\"],\n  [\"314\u00d73=942\", \"503\u00d75=2515\"],\n  [\"607\u00d78=4856\", \"294\u00d77=2058\"],\n  [\"505\u00d78=4040\", \"883\u00d78=7064\"],\n  [\"202\u00d72=404\", \"160\u00d72=320\"],\n  [\"503\u00d77=3521\", \"511\u00d76=3066\"],\n  [\"816\u00d74=3264\", \"656\u00d78=5248\"],\n  [\"722\u00d73=2166\", \"335\u00d77=2345\"],\n  [\"947\u00d76=5682\", \"498\u00d72=996\"],\n  [\"925\u00d76=5550\", \"878\u00d76=5268\"],\n  [\"160\u00d74=640\", \"272\u00d75=1360\"],\n  [\"342\u00d72=684\", \"145\u00d78=1160\"],\n  [\"681\u00d75=3405\", \"743\u00d78=5944\"],\n  [\"780\u00d74=3120\", \"815\u00d74=3260\"],\n  [\"962\u00d72=1924\", \"954\u00d72=1908\"],\n  [\"944\u00d72=1888\", \"155\u00d79=1395\"],\n  [\"307\u00d76=1842\", \"729\u00d76=4374\"],\n  [\"583\u00d76=3498\", \"605\u00d79=5445\"],\n  [\"109\u00d79=981\", \"467\u00d77=3269\"],\n  [\"743\u00d77=5201\", \"497\u00d78=3976\"],\n  [\"118\u00d77=826\", \"722\u00d78=5776\"],\n  [\"797\u00d75=3985\", \"343\u00d72=686\"],\n  [\"322\u00d72=644\", \"227\u00d74=908\"],\n  [\"306\u00d79=2754\", \"455\u00d75=2275\"],\n];\n\nconst body = context.document.body;\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(after, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"393\u00d73=1179\", \"297\u00d77=2079\"),\n    @(\"693\u00d72=1386\", \"402\u00d78=3216\"),\n    @(\"314\u00d73=942\", \"503\u00d75=2515\"),\n    @(\"607\u00d78=4856\", \"294\u00d77=2058\"),\n    @(\"505\u00d78=4040\", \"883\u00d78=7064\"),\n    @(\"202\u00d72=404\", \"160\u00d72=320\"),\n    @(\"503\u00d77=3521\", \"511\u00d76=3066\"),\n    @(\"816\u00d74=3264\", \"656\u00d78=5248\"),\n    @(\"722\u00d73=2166\", \"335\u00d77=2345\"),\n    @(\"947\u00d76=5682\", \"498\u00d72=996\"),\n    @(\"925\u00d76=5550\", \"878\u00d76=5268\"),\n    @(\"160\u00d74=640\", \"272\u00d75=1360\"),\n    @(\"342\u00d72=684\", \"145\u00d78=1160\"),\n    @(\"681\u00d75=3405\", \"743\u00d78=5944\"),\n    @(\"780\u00d74=3120\", \"815\u00d74=3260\"),\n    @(\"962\u00d72=1924\", \"954\u00d72=1908\"),\n    @(\"944\u00d72=1888\", \"155\u00d79=1395\"),\n    @(\"307\u00d76=1842\", \"729\u00d76=4374\"),\n    @(\"583\u00d76=3498\", \"605\u00d79=5445\"),\n    @(\"109\u00d79=981\", \"467\u00d77=3269\"),\n    @(\"743\u00d77=5201\", \"497\u00d78=3976\"),\n    @(\"118\u00d77=826\", \"722\u00d78=5776\"),\n    @(\"797\u00d75=3985\", \"343\u00d72=686\"),\n    @(\"322\u00d72=644\", \"227\u00d74=908\"),\n    @(\"306\u00d79=2754\", \"455\u00d75=2275\"),\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}"}
